# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, and fixes the OKB / Avalanche row ordering swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ D='39.085.06'; E='  -4.62%  ' }
    3 = @{ D='2.234.29'; E='  -7.23%  ' }
    4 = @{ D='1.00'; E='  -0.06%  ' }
    5 = @{ D='296.19'; E='  -5.87%  ' }
    6 = @{ D='80.45'; E='  -8.75%  ' }
    7 = @{ D='0.508'; E='  -4.85%  ' }
    8 = @{ E='  -0.02%  ' }
    9 = @{ E='  -7.48%  ' }
    10 = @{ D='0.0774'; E='  -6.82%  ' }
    11 = @{ B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='27.90'; E='  -10.55%  ' }
    12 = @{ B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='46.45'; E='  -12.90%  ' }
    13 = @{ D='0.107'; E='  -1.37%  ' }
    14 = @{ D='2.582.93'; E='  -7.12%  ' }
    15 = @{ E='  -9.90%  ' }
    16 = @{ E='  -8.96%  ' }
    17 = @{ D='2.240.29'; E='  -7.09%  ' }
    18 = @{ D='0.713'; E='  -7.07%  ' }
    19 = @{ D='38.965.55'; E='  -4.58%  ' }
    20 = @{ D='0.0₃0857'; E='  -6.65%  ' }
    22 = @{ D='65.39'; E='  -7.16%  ' }
    23 = @{ D='9.90'; E='  -8.49%  ' }
    24 = @{ D='226.51'; E='  -5.23%  ' }
    25 = @{ E='  -0.01%  ' }
    26 = @{ E='  -10.16%  ' }
    27 = @{ E='  -6.37%  ' }
    28 = @{ D='22.22'; E='  -7.01%  ' }
    29 = @{ E='  -1.99%  ' }
    30 = @{ D='8.84'; E='  -6.69%  ' }
    31 = @{ D='147.56'; E='  -5.95%  ' }
    32 = @{ D='31.39'; E='  -7.88%  ' }
    33 = @{ D='1.00'; E='  -0.19%  ' }
    34 = @{ D='4.76'; E='  -9.35%  ' }
    35 = @{ D='2.30'; E='  -6.58%  ' }
    36 = @{ E='  -7.33%  ' }
    37 = @{ E='  -3.92%  ' }
    38 = @{ D='2.62'; E='  -7.97%  ' }
    39 = @{ D='0.0944'; E='  -4.86%  ' }
    40 = @{ D='14.64'; E='  -8.97%  ' }
    41 = @{ E='  -9.06%  ' }
    42 = @{ D='3.63'; E='  -5.70%  ' }
    43 = @{ D='1.911.43'; E='  -3.77%  ' }
    44 = @{ D='2.19'; E='  -4.17%  ' }
    45 = @{ D='0.0253'; E='  -7.02%  ' }
    46 = @{ D='16.22'; E='  -9.57%  ' }
    47 = @{ D='8.93'; E='  -4.04%  ' }
    48 = @{ E='  -11.35%  ' }
    49 = @{ D='2.464.73'; E='  -6.92%  ' }
    50 = @{ D='87.58'; E='  -6.42%  ' }
    51 = @{ D='66.31'; E='  -10.09%  ' }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]

    if ($cols.ContainsKey('B')) {
        $ws.Range("B$row").Value = $cols['B']
    }
    if ($cols.ContainsKey('C')) {
        $ws.Range("C$row").Value = $cols['C']
    }
    if ($cols.ContainsKey('D')) {
        # Price column holds numeric-looking text (e.g. "1.00", "296.19").
        # Force text format first so Excel doesn't silently coerce it to a
        # number and drop the original formatting.
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = "@"
        $dCell.Value = $cols['D']
        $dCell.Style = "Normal"
    }
    if ($cols.ContainsKey('E')) {
        $ws.Range("E$row").Value = $cols['E']
    }
}
